$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the end time on row 22 (18:50:45 -> 18:50:00), stored as a day fraction
$ws.Range("D22").Value = 0.784722222222222

# Row 29 previously had no start/end time recorded; fill in start (16:30) and end (18:45) times
$ws.Range("C29").Value = 0.6875
$ws.Range("D29").Value = 0.78125

# Update the selected/active cell to reflect where the author was last editing
$ws.Range("C5").Select()
